$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The original sheet had a header row (row 1) and a single data row (row 2,
# "Acacia Breeze"). Two new projects are inserted *before* "Acacia Breeze"
# (new rows 2 and 3), and one new project is appended *after* it (new row 5).
# ---------------------------------------------------------------------------

# Insert two blank rows at 2:3 - this pushes "Acacia Breeze" from row 2 to row 4.
$ws.Range("A2:A3").EntireRow.Insert()

# New row 2: Bukit Merah Ridges
$ws.Range("A2").Value2 = "Bukit Merah Ridges"
$ws.Range("B2").Value2 = "Bukit_Merah"
$ws.Range("C2").Value2 = "2-Room"
$ws.Range("D2").Value2 = 3
$ws.Range("E2").Value2 = 400000
$ws.Range("F2").Value2 = "3-Room"
$ws.Range("G2").Value2 = 5
$ws.Range("H2").Value2 = 650000
$ws.Range("K2").Value2 = "Michael"
$ws.Range("L2").Value2 = 2
$ws.Range("M2").Value2 = "David"

# New row 3: Bukit Panjang Heights
$ws.Range("A3").Value2 = "Bukit Panjang Heights"
$ws.Range("B3").Value2 = "Bukit_Panjang"
$ws.Range("C3").Value2 = "2-Room"
$ws.Range("D3").Value2 = 1
$ws.Range("E3").Value2 = 350000
$ws.Range("F3").Value2 = "3-Room"
$ws.Range("G3").Value2 = 1
$ws.Range("H3").Value2 = 460000
$ws.Range("K3").Value2 = "Jessica"
$ws.Range("L3").Value2 = 10
$ws.Range("M3").Value2 = "Emily"

# Row 4 keeps the pre-existing "Acacia Breeze" data (shifted down automatically).

# Append a new row 5 after "Acacia Breeze": Nanyang Gardens
$ws.Range("A5:M5").EntireRow.Insert()
$ws.Range("A5").Value2 = "Nanyang Gardens"
$ws.Range("B5").Value2 = "Jurong_West"
$ws.Range("C5").Value2 = "3-Room"
$ws.Range("D5").Value2 = 2
$ws.Range("E5").Value2 = 125000
$ws.Range("F5").Value2 = "2-Room"
$ws.Range("G5").Value2 = 1
$ws.Range("H5").Value2 = 400000
$ws.Range("K5").Value2 = "Michael"
$ws.Range("L5").Value2 = 10
$ws.Range("M5").Value2 = "Daniel, Emily, David"

# Date columns (I = opening date, J = closing date) need to keep the existing
# short-date style (style index 1 / numFmtId 14) rather than minting a new
# number format, so copy formats from the still-correctly-styled row 4 dates
# before writing the new date values.
$ws.Range("I4:J4").Copy()
$ws.Range("I2:J2").PasteSpecial(-4122)
$ws.Range("I3:J3").PasteSpecial(-4122)
$ws.Range("I5:J5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I2").Value2 = 45599
$ws.Range("J2").Value2 = 45613
$ws.Range("I3").Value2 = 45599
$ws.Range("J3").Value2 = 45613
$ws.Range("I5").Value2 = 45757
$ws.Range("J5").Value2 = 45777

# Column M widens to fit the longest Officer list ("Daniel, Emily, David").
$ws.Range("M1").ColumnWidth = 15.666666666666666

# Final cursor position left on F6.
$ws.Range("F6").Select() | Out-Null
